$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns I and J (COORD_X/COORD_Y) for every data row,
# inverting the coordinate columns as described in the commit message.
for ($r = 2; $r -le 22; $r++) {
    $cellI = $ws.Cells.Item($r, 9)
    $cellJ = $ws.Cells.Item($r, 10)
    $i = $cellI.Value2
    $j = $cellJ.Value2
    $cellI.Value2 = $j
    $cellJ.Value2 = $i
}

# Give the (now populated) columns I and J explicit widths.
$ws.Columns.Item(9).ColumnWidth = 15.66
$ws.Columns.Item(10).ColumnWidth = 18.0

# Move the active selection to H23.
$ws.Range("H23").Select() | Out-Null
